$wb = $excel.ActiveWorkbook
$origSheet = $wb.ActiveSheet
$ws = $wb.Worksheets.Item(5)
$ws.Activate()

# --- Table 2 (rows 9-13) first, matching the order the author typed them in ---
$ws.Range("E9").Value = " Time1"
$ws.Range("F9").Value = "Time2"
$ws.Range("G9").Value = "Time3"
$ws.Range("H9").Value = "Average "

$ws.Range("F10").Value = 0.0390625
$ws.Range("G10").Value = 0.03125
$ws.Range("H10").Formula = "=AVERAGE(E10:G10)"

$ws.Range("F11").Value = 0.0078125
$ws.Range("G11").Value = 0.0078125
$ws.Range("H11").Formula = "=AVERAGE(E11:G11)"

$ws.Range("F12").Value = 0.0078125
$ws.Range("G12").Value = 0.0078125
$ws.Range("H12").Formula = "=AVERAGE(E12:G12)"

$ws.Range("F13").Value = 0.0234375
$ws.Range("G13").Value = 0.0234375
$ws.Range("H13").Formula = "=AVERAGE(E13:G13)"

# --- Table 1 (rows 2-6) ---
$ws.Range("L2").Value = " Time1"
$ws.Range("M2").Value = "Time2"
$ws.Range("N2").Value = "Time3"
$ws.Range("O2").Value = "Average"

$ws.Range("M3").Value = 0.0546875
$ws.Range("N3").Value = 0.0390625
$ws.Range("O3").Formula = "=AVERAGE(L3:N3)"

$ws.Range("M4").Value = 0.0234375
$ws.Range("N4").Value = 0.015625
$ws.Range("O4").Formula = "=AVERAGE(L4:N4)"

$ws.Range("M5").Value = 0.0078125
$ws.Range("N5").Value = 0.015625
$ws.Range("O5").Formula = "=AVERAGE(L5:N5)"

$ws.Range("M6").Value = 0.03125
$ws.Range("N6").Value = 0.0390625
$ws.Range("O6").Formula = "=AVERAGE(L6:N6)"

$ws.Range("N9").Select()
$origSheet.Activate()
